$wb = $excel.ActiveWorkbook

# --- A20: update the current selection (active cell) ---
$wsA20 = $wb.Worksheets.Item("A20")
$wsA20.Range("G18").Select()

# --- B1: new stage-3 rows (shared strings 88, 89), resized column, new selection ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "721FECB35598A390E92659848A7C7D55BAE2504748173EB02E22CE678B877E76"
$wsB1.Range("A3").Value = "E97E1C1F81926FB07A7FEEB3483BE58CC4456A4468DCA1DFFB965F6552B43091"
$wsB1.Columns.Item(1).ColumnWidth = 71.25
$wsB1.Range("A4").Select()

# --- B2: new stage-3 rows (shared strings 90, 91 -- A3 written first), resized column, new selection ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A3").Value = "FBAB86A707514C51B7EC6206A36C915ECE68DB9A47F6A89E24FDA1B1AC2D634A"
$wsB2.Range("A2").Value = "B10DE0190BCB0D4779CA0D03CF4BEA6FD20ED34B2EE284612823F78A02D9E6CF"
$wsB2.Columns.Item(1).ColumnWidth = 72.59
$wsB2.Range("A2").Select()

# --- Re-activate B1 last so it ends up as the selected/active tab ---
$wsB1.Activate()
$wsB1.Range("A4").Select()
